# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H,I,J,K,L,M,N) across the per-job Leve tables. Columns A-G (leve/item info)
# are untouched.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2509.4707
$ws.Range("L17").Value = 7865.0772
$ws.Range("J17").Value = 2621.6924
$ws.Range("N17").Value = -8201.0772
$ws.Range("H40").Value = 5147.0835
$ws.Range("L40").Value = 6416.6
$ws.Range("N40").Value = -6766.6
$ws.Range("J40").Value = 6416.6
$ws.Range("H132").Value = 4034.318
$ws.Range("M132").Value = -10671.6581
$ws.Range("K132").Value = 13201.6581
$ws.Range("I132").Value = 4400.5527
$ws.Range("H135").Value = 796.8
$ws.Range("M135").Value = -4741.5
$ws.Range("K135").Value = 7276.5
$ws.Range("I135").Value = 808.5
$ws.Range("I137").Value = 1962.8077
$ws.Range("M137").Value = -3338.4231
$ws.Range("L137").Value = 15174.4995
$ws.Range("K137").Value = 5888.4231
$ws.Range("N137").Value = -20274.4995
$ws.Range("J137").Value = 5058.1665
$ws.Range("H137").Value = 2543.1875
$ws.Range("L138").Value = 1066320.84
$ws.Range("M138").Value = -5975
$ws.Range("N138").Value = -1076600.84
$ws.Range("J138").Value = 355440.28
$ws.Range("K138").Value = 11115
$ws.Range("H138").Value = 253437.05
$ws.Range("I138").Value = 3705

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 9797
$ws.Range("K13").Value = 9797
$ws.Range("M13").Value = -9653
$ws.Range("I13").Value = 9797
$ws.Range("M32").Value = -1484.1052
$ws.Range("H32").Value = 2382.52
$ws.Range("K32").Value = 1771.1052
$ws.Range("I32").Value = 1771.1052
$ws.Range("N43").Value = -80620.5
$ws.Range("L43").Value = 79994.5
$ws.Range("H43").Value = 83329.664
$ws.Range("J43").Value = 79994.5
$ws.Range("L44").Value = 77000
$ws.Range("H44").Value = 77000
$ws.Range("N44").Value = -77976
$ws.Range("J44").Value = 77000
$ws.Range("N80").Value = -86996
$ws.Range("H80").Value = 85000
$ws.Range("L80").Value = 85000
$ws.Range("J80").Value = 85000
$ws.Range("J83").Value = 85000
$ws.Range("N83").Value = -264984
$ws.Range("H83").Value = 85000
$ws.Range("L83").Value = 255000
$ws.Range("I122").Value = 3905.889
$ws.Range("K122").Value = 11717.667
$ws.Range("M122").Value = -9267.667000000001
$ws.Range("H122").Value = 4629.516
$ws.Range("H132").Value = 2627.84
$ws.Range("M132").Value = -4361.9171
$ws.Range("K132").Value = 6891.9171
$ws.Range("I132").Value = 2297.3057

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K105").Value = 529639.2
$ws.Range("M105").Value = -527892.2
$ws.Range("L105").Value = 31253796
$ws.Range("J105").Value = 31253796
$ws.Range("N105").Value = -31257290
$ws.Range("H105").Value = 9633093
$ws.Range("I105").Value = 529639.2
$ws.Range("K134").Value = 9691.2855
$ws.Range("I134").Value = 3230.4285
$ws.Range("H134").Value = 3540.5833
$ws.Range("M134").Value = -7156.2855

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 2398.9048
$ws.Range("N16").Value = -3779.6
$ws.Range("H16").Value = 2659.1292
$ws.Range("K16").Value = 2398.9048
$ws.Range("M16").Value = -2111.9048
$ws.Range("J16").Value = 3205.6
$ws.Range("L16").Value = 3205.6
$ws.Range("H31").Value = 3097.2976
$ws.Range("L31").Value = 4722.4
$ws.Range("N31").Value = -5312.4
$ws.Range("J31").Value = 4722.4
$ws.Range("N34").Value = -5126.4
$ws.Range("L34").Value = 4722.4
$ws.Range("J34").Value = 4722.4
$ws.Range("H34").Value = 3097.2976
$ws.Range("I58").Value = 1857.7333
$ws.Range("K58").Value = 1857.7333
$ws.Range("M58").Value = -1654.7333
$ws.Range("H58").Value = 3084.96
$ws.Range("K86").Value = 2283.75
$ws.Range("H86").Value = 3024
$ws.Range("M86").Value = -1160.75
$ws.Range("I86").Value = 2283.75
$ws.Range("I89").Value = 2283.75
$ws.Range("M89").Value = -5802.75
$ws.Range("K89").Value = 11418.75
$ws.Range("H89").Value = 3024
$ws.Range("K113").Value = 2398.9048
$ws.Range("N113").Value = -7545.6
$ws.Range("L113").Value = 3205.6
$ws.Range("J113").Value = 3205.6
$ws.Range("I113").Value = 2398.9048
$ws.Range("M113").Value = -228.9047999999998
$ws.Range("H113").Value = 2659.1292
$ws.Range("K134").Value = 5218.028700000001
$ws.Range("I134").Value = 1739.3429
$ws.Range("H134").Value = 2321.7917
$ws.Range("M134").Value = -2683.028700000001
$ws.Range("K136").Value = 5573.199900000001
$ws.Range("M136").Value = -3023.199900000001
$ws.Range("H136").Value = 3084.96
$ws.Range("I136").Value = 1857.7333
$ws.Range("L141").Value = 382491.84
$ws.Range("H141").Value = 382491.84
$ws.Range("J141").Value = 382491.84
$ws.Range("N141").Value = -392851.84

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M44").Value = -2564.5
$ws.Range("I44").Value = 987.5
$ws.Range("K44").Value = 2962.5
$ws.Range("H44").Value = 1565.3572
$ws.Range("N97").Value = -152333.75
$ws.Range("L97").Value = 151341.75
$ws.Range("H97").Value = 743254.9
$ws.Range("J97").Value = 50447.25
$ws.Range("N122").Value = -22660.0006
$ws.Range("L122").Value = 17760.0006
$ws.Range("J122").Value = 1973.3334
$ws.Range("H122").Value = 1780
$ws.Range("N125").Value = -36840
$ws.Range("J125").Value = 9000
$ws.Range("L125").Value = 27000
$ws.Range("H125").Value = 9000
$ws.Range("K128").Value = 396805.26
$ws.Range("H128").Value = 132268.42
$ws.Range("I128").Value = 132268.42
$ws.Range("M128").Value = -391825.26
$ws.Range("I131").Value = 12501.77
$ws.Range("H131").Value = 5248.6904
$ws.Range("K131").Value = 37505.31
$ws.Range("M131").Value = -32465.31

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("J53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("J57").Value = 87653
$ws.Range("N57").Value = -89293
$ws.Range("H57").Value = 87653
$ws.Range("L57").Value = 87653
$ws.Range("I102").Value = 1747.421
$ws.Range("K102").Value = 1747.421
$ws.Range("J102").Value = 21480.8
$ws.Range("H102").Value = 5858.5415
$ws.Range("L102").Value = 21480.8
$ws.Range("M102").Value = -125.421
$ws.Range("N102").Value = -24724.8

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M40").Value = -36567.516
$ws.Range("H40").Value = 22011.25
$ws.Range("K40").Value = 36703.516
$ws.Range("L40").Value = 9701.513999999999
$ws.Range("I40").Value = 36703.516
$ws.Range("N40").Value = -9973.513999999999
$ws.Range("J40").Value = 9701.513999999999
$ws.Range("I100").Value = 3911.423
$ws.Range("J100").Value = 6732.3335
$ws.Range("K100").Value = 3911.423
$ws.Range("M100").Value = -3370.423
$ws.Range("H100").Value = 4203.241
$ws.Range("L100").Value = 6732.3335
$ws.Range("N100").Value = -7814.3335
$ws.Range("N122").Value = -20204.5
$ws.Range("I122").Value = 7751
$ws.Range("L122").Value = 15304.5
$ws.Range("K122").Value = 23253
$ws.Range("J122").Value = 5101.5
$ws.Range("M122").Value = -20803
$ws.Range("H122").Value = 5984.6665
$ws.Range("K136").Value = 16157.4
$ws.Range("M136").Value = -13607.4
$ws.Range("H136").Value = 6838.85
$ws.Range("I136").Value = 5385.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J107").Value = 949.5
$ws.Range("L107").Value = 2848.5
$ws.Range("I107").Value = 591.4
$ws.Range("N107").Value = -6688.5
$ws.Range("M107").Value = 145.8000000000002
$ws.Range("H107").Value = 651.0833
$ws.Range("K107").Value = 1774.2
$ws.Range("N122").Value = -250018900
$ws.Range("I122").Value = 2271.389
$ws.Range("L122").Value = 250014000
$ws.Range("K122").Value = 6814.167
$ws.Range("J122").Value = 83338000
$ws.Range("M122").Value = -4364.167
$ws.Range("H122").Value = 11907376
$ws.Range("H126").Value = 1801.68
$ws.Range("K126").Value = 5230.6362
$ws.Range("M126").Value = -2760.6362
$ws.Range("L126").Value = 6684
$ws.Range("J126").Value = 2228
$ws.Range("I126").Value = 1743.5454
$ws.Range("N126").Value = -11624
$ws.Range("H132").Value = 4897.6387
$ws.Range("M132").Value = -11933.723
$ws.Range("K132").Value = 14463.723
$ws.Range("I132").Value = 4821.241
$ws.Range("N136").Value = -11265
$ws.Range("J136").Value = 2055
$ws.Range("H136").Value = 40001660
$ws.Range("L136").Value = 6165
